# dropping health notes as a field
#
# The "Health notes" header lives in column O (column 15) of the single
# "Blank template" worksheet. Deleting the whole column removes it from
# the shared-string table and shifts every later field (Host life stage,
# Age, ... GenBank accession) one column to the left, exactly like using
# Excel's "Delete Column" on O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(15).Delete()

# The sheet's used range originally extended through AO1 (41 columns).
# After the shift, the last real header ("GenBank accession") now sits in
# AN1, so AO1 becomes a trailing blank cell - touch its number format so
# the cell (and therefore the A1:AO1 extent) is preserved instead of being
# dropped.
$ws.Range("AO1").NumberFormat = "General"

# Reflect the edited column in the current selection.
[void]$ws.Range("O1").Select()
